# Cross_sell.xlsx - "customer overview feature file added. made few changes to cross sell page tests."
#
# 1. Cross_sell_Filter (sheet 1): swap the "Segmentation"/"Date Range" picks between
#    row 2 and row 3 (B2/D2 <-> B3/D3), and make the second sheet the active tab/view.
# 2. Cross_Sell_Report_EtoE (sheet 2): insert a new "Selected Drivers" column (D),
#    populate it, swap rows 4 and 7's data back into the "90%-100%"/"50%-100%" order,
#    and set the new active selection / page setup.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Cross_sell_Filter
$ws2 = $wb.Worksheets.Item(2)   # Cross_Sell_Report_EtoE

# --- Sheet 1: Cross_sell_Filter ---------------------------------------------
# Move the "Overall" segmentation pick (with its formatting) from B3 up to B2,
# clear B3's value but keep its (now shared) look.
$ws1.Range("B3").Copy($ws1.Range("B2"))
$ws1.Range("B3").ClearContents()

# Swap the "Date Range" picks between row 2 and row 3.
$ws1.Range("D2").Value = "Within 180 days"
$ws1.Range("D3").Value = "Within 90 days"

$ws1.Range("B8").Select()

# --- Sheet 2: Cross_Sell_Report_EtoE ----------------------------------------
# Insert a new column D ("Selected Drivers") between "Customer Probability" and
# "Drivers Title".
$ws2.Columns.Item(4).Insert()
$ws2.Columns.Item(4).ColumnWidth = 29.7

$ws2.Range("D1").Value = "Selected Drivers"
$ws2.Range("D3").Value = "NA"
$ws2.Range("D2").Value = "Driver 1,Driver 2,Driver 3,Driver 4,Driver 5"
$ws2.Range("D4").Value = "Driver 1,Driver 2,Driver 3,Driver 4,Driver 5"
$ws2.Range("D5").Value = "NA"
$ws2.Range("D6").Value = "NA"
$ws2.Range("D7").Value = "Driver 1,Driver 2,Driver 3,Driver 4,Driver 5"

# Rows 4 and 7 swap their Customer-Probability-Range / Customer-Probability /
# Drivers-Title / Profile-Title content (the "100%-50%" combo row moves from 4
# down to 7, and the "100%-90%" combo row moves from 7 up to 4).
$ws2.Range("B4").Value = "100%-90%"
$ws2.Range("C4").Value = "Extreme Likely"
$ws2.Range("E4").Value = "Cross-Sell Drivers for Top 90% to 100%"
$ws2.Range("F4").Value = "Profiles for Top 90% to 100%"

$ws2.Range("B7").Value = "100%-50%"
$ws2.Range("C7").Value = "Extreme Likely,High Likely,Likely"
$ws2.Range("E7").Value = "Cross-Sell Drivers for Top 50% to 100%"
$ws2.Range("F7").Value = "Profiles for Top 50% to 100%"

$ws2.Range("D12").Select()

# Sheet 2 becomes the active tab/view.
$ws2.Activate()

# New print setup for the report sheet.
$ws2.PageSetup.Orientation = 1
